# Adds the "HashSet (ex. 8.3)" subsection at the end of the "Collections"
# section, mirroring the structure of the existing "List for whole numbers
# (ex. 8.1)" / "PhoneBook (ex. 8.2)" subsections: a Heading2 title paragraph
# followed by a FirstParagraph body paragraph, the pair wrapped in a bookmark
# named "hashset-ex.-8.3".

$d = $word.ActiveDocument

# New content is appended right after the last paragraph currently in the
# document ("/exit", end of the PhoneBook bullet list), i.e. at document end.
$endRange = $d.Range($d.Content.End, $d.Content.End)

# --- Heading2 paragraph: "HashSet (ex. 8.3)" ---------------------------
$headingPara = $d.Paragraphs.Add($endRange)
$headingPara.Range.Text = "HashSet (ex. 8.3)"
$headingPara.Style = "Heading 2"

$bookmarkStartPos = $headingPara.Range.Start

# --- FirstParagraph paragraph: body text --------------------------------
$bodyRange = $d.Range($d.Content.End, $d.Content.End)
$bodyPara = $d.Paragraphs.Add($bodyRange)
$bodyPara.Range.Text = "Realisation just due task requirements."
$bodyPara.Style = "First Paragraph"

$bookmarkEndPos = $bodyPara.Range.End

# Wrap both new paragraphs in a single bookmark, the same way the sibling
# subsections above (list-for-whole-numbers-ex.-8.1, phonebook-ex.-8.2) are
# each wrapped in their own bookmark.
$newSectionRange = $d.Range($bookmarkStartPos, $bookmarkEndPos)
$d.Bookmarks.Add("hashset-ex.-8.3", $newSectionRange)
